# Update cryptos list: refresh Price (column D) and Volume(1h) (column E) values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = "33.999.37";   E = "  +0.63%  " },
    @{ Row = 3;  D = "1.777.38";    E = "  -1.45%  " },
    @{ Row = 4;  D = $null;         E = "  +0.26%  " },
    @{ Row = 5;  D = "226.01";      E = "  -0.53%  " },
    @{ Row = 6;  D = "0.547";       E = "  +1.80%  " },
    @{ Row = 7;  D = $null;         E = "  +0.27%  " },
    @{ Row = 8;  D = "31.15";       E = "  +0.84%  " },
    @{ Row = 9;  D = "46.51";       E = "  -1.54%  " },
    @{ Row = 10; D = $null;         E = "  +0.55%  " },
    @{ Row = 11; D = "0.0656";      E = "  -0.70%  " },
    @{ Row = 12; D = "0.0930";      E = "  +0.21%  " },
    @{ Row = 13; D = "2.037.74";    E = "  -1.34%  " },
    @{ Row = 14; D = $null;         E = "  +11.11%  " },
    @{ Row = 15; D = "1.790.40";    E = "  -0.78%  " },
    @{ Row = 16; D = "0.625";       E = "  -0.97%  " },
    @{ Row = 17; D = "33.997.62";   E = "  +0.79%  " },
    @{ Row = 18; D = "4.20";        E = "  -0.50%  " },
    @{ Row = 19; D = "68.94";       E = $null },
    @{ Row = 20; D = "251.83";      E = "  -1.00%  " },
    @{ Row = 21; D = "0.0₃0736";    E = "  -0.30%  " },
    @{ Row = 22; D = $null;         E = "  +0.32%  " },
    @{ Row = 23; D = "10.36";       E = "  -0.24%  " },
    @{ Row = 24; D = "4.20";        E = "  -2.10%  " },
    @{ Row = 25; D = $null;         E = "  -1.88%  " },
    @{ Row = 26; D = "156.12";      E = "  -0.98%  " },
    @{ Row = 27; D = "16.38";       E = "  +0.17%  " },
    @{ Row = 28; D = "6.98";        E = "  -0.35%  " },
    @{ Row = 29; D = "0.113";       E = "  -0.65%  " },
    @{ Row = 30; D = $null;         E = "  +0.28%  " },
    @{ Row = 31; D = "3.75";        E = "  -1.05%  " },
    @{ Row = 32; D = "0.0513";      E = "  +1.27%  " },
    @{ Row = 33; D = $null;         E = "  +0.59%  " },
    @{ Row = 34; D = "3.56";        E = "  +2.36%  " },
    @{ Row = 35; D = $null;         E = "  +2.49%  " },
    @{ Row = 36; D = "1.444.52";    E = "  -5.54%  " },
    @{ Row = 37; D = $null;         E = "  -0.76%  " },
    @{ Row = 38; D = "0.626";       E = "  +2.32%  " },
    @{ Row = 39; D = $null;         E = "  +1.32%  " },
    @{ Row = 40; D = "2.85";        E = "  +2.09%  " },
    @{ Row = 41; D = "82.40";       E = "  -0.93%  " },
    @{ Row = 42; D = "2.35";        E = $null },
    @{ Row = 43; D = "0.889";       E = "  -1.01%  " },
    @{ Row = 44; D = "2.05";        E = "  -2.37%  " },
    @{ Row = 45; D = $null;         E = "  -1.93%  " },
    @{ Row = 46; D = $null;         E = "  -0.74%  " },
    @{ Row = 47; D = "1.934.66";    E = "  -1.00%  " },
    @{ Row = 48; D = $null;         E = "  +2.68%  " },
    @{ Row = 49; D = $null;         E = "  +0.39%  " },
    @{ Row = 50; D = "11.72";       E = "  +5.33%  " },
    @{ Row = 51; D = "49.91";       E = "  -3.61%  " }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($null -ne $u.D) {
        $cell = $ws.Cells.Item($r, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.Style = "Normal"
    }
    if ($null -ne $u.E) {
        $cell = $ws.Cells.Item($r, 5)
        $cell.NumberFormat = "@"
        $cell.Value = $u.E
        $cell.Style = "Normal"
    }
}
